# Generate Report for Handback
#
# The handback-status report records, for each localized language sheet,
# the "Correspond Handoff Datetime" and "Correspond Handback DateTime" for
# every handed-back file. Re-running the report generator refreshed the
# timestamps for the file "35efe67b-7699-461d-923e-8925f6541628" (row 2 of
# each language sheet) on both the "zh-cn" and "de-de" sheets, while the
# "86b20423-8bb9-4abc-9c66-0eed75ac83da" row (row 3) stayed in sync with its
# previous handback and kept its existing timestamps.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-20 18:36:24"
$wsZhCn.Range("H2").Value = "2016-03-20 18:36:50"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-20 18:36:27"
$wsDeDe.Range("H2").Value = "2016-03-20 18:36:55"
